# Update simulation output values in the "Comparison" sheet to reflect the
# latest projection run (headcount/eligibility/participation/contribution figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9465
$ws.Range("E2").Value = 8368
$ws.Range("F2").Value = 0.8840993132593766
$ws.Range("G2").Value = 0.8828866849546317
$ws.Range("H2").Value = 0.0958304262959481
$ws.Range("I2").Value = 0.08460740739021877
$ws.Range("J2").Value = 40617149.71070025
$ws.Range("K2").Value = 14113742.80200613
$ws.Range("M2").Value = 14113742.80200613
$ws.Range("N2").Value = 54730892.51270638
$ws.Range("O2").Value = 800741935.3172001
$ws.Range("P2").Value = 783042128.3132001
$ws.Range("Q2").Value = 0.01762583196846711
$ws.Range("R2").Value = 0.01802424453510492

# Row 3
$ws.Range("D3").Value = 9642
$ws.Range("E3").Value = 8557
$ws.Range("F3").Value = 0.8874714789462768
$ws.Range("G3").Value = 0.8852679495137595
$ws.Range("H3").Value = 0.09415192489269762
$ws.Range("I3").Value = 0.08334968149253193
$ws.Range("J3").Value = 42347526.09196078
$ws.Range("K3").Value = 14709013.98571959
$ws.Range("M3").Value = 14709013.98571959
$ws.Range("N3").Value = 57056540.07768038
$ws.Range("O3").Value = 836091164.5578281
$ws.Range("P3").Value = 818610988.533758
$ws.Range("Q3").Value = 0.01759259589054328
$ws.Range("R3").Value = 0.0179682586622315

# Row 4
$ws.Range("D4").Value = 9836
$ws.Range("E4").Value = 8714
$ws.Range("F4").Value = 0.8859292395282635
$ws.Range("G4").Value = 0.8839521201054981
$ws.Range("H4").Value = 0.09305289739007716
$ws.Range("I4").Value = 0.08225430592991807
$ws.Range("J4").Value = 44259719.44997451
$ws.Range("K4").Value = 15345188.53406264
$ws.Range("M4").Value = 15345188.53406264
$ws.Range("N4").Value = 59604907.98403715
$ws.Range("O4").Value = 874653871.4285319
$ws.Range("P4").Value = 857204923.4225781
$ws.Range("Q4").Value = 0.01754429841944225
$ws.Range("R4").Value = 0.017901423702507

# Row 5
$ws.Range("D5").Value = 10027
$ws.Range("E5").Value = 8908
$ws.Range("F5").Value = 0.8884013164455968
$ws.Range("G5").Value = 0.8860155162124528
$ws.Range("H5").Value = 0.09171498176794257
$ws.Range("I5").Value = 0.08126089691553934
$ws.Range("J5").Value = 46251536.29016398
$ws.Range("K5").Value = 16008023.2459684
$ws.Range("M5").Value = 16008023.2459684
$ws.Range("N5").Value = 62259559.53613237
$ws.Range("O5").Value = 913201104.5186434
$ws.Range("P5").Value = 895714999.0556703
$ws.Range("Q5").Value = 0.017529570613481
$ws.Range("R5").Value = 0.017871782054387

# Row 6
$ws.Range("D6").Value = 10222
$ws.Range("E6").Value = 9062
$ws.Range("F6").Value = 0.8865192721580903
$ws.Range("G6").Value = 0.8837526818802418
$ws.Range("H6").Value = 0.09064143435825292
$ws.Range("I6").Value = 0.08010461070357792
$ws.Range("J6").Value = 48302938.50549269
$ws.Range("K6").Value = 16662679.69380575
$ws.Range("M6").Value = 16662679.69380575
$ws.Range("N6").Value = 64965618.19929844
$ws.Range("O6").Value = 954323256.6988841
$ws.Range("P6").Value = 936731430.8257025
$ws.Range("Q6").Value = 0.01746020499536385
$ws.Range("R6").Value = 0.01778810782415838

# Row 7
$ws.Range("D7").Value = 9453
$ws.Range("E7").Value = 8379
$ws.Range("F7").Value = 0.8863852745160267
$ws.Range("G7").Value = 0.8840472673559823
$ws.Range("H7").Value = 0.09675112722510387
$ws.Range("I7").Value = 0.08553256963696407
$ws.Range("J7").Value = 41181848.30283703
$ws.Range("K7").Value = 14396092.09807451
$ws.Range("M7").Value = 14396092.09807451
$ws.Range("N7").Value = 55577940.40091153
$ws.Range("O7").Value = 800091928.4872
$ws.Range("P7").Value = 782392121.4832001
$ws.Range("Q7").Value = 0.01799304753054364
$ws.Range("R7").Value = 0.01840009849636968

# Row 8
$ws.Range("D8").Value = 9644
$ws.Range("E8").Value = 8563
$ws.Range("F8").Value = 0.887909581086686
$ws.Range("G8").Value = 0.8858886819780675
$ws.Range("H8").Value = 0.09538540115674442
$ws.Range("I8").Value = 0.08450084731069751
$ws.Range("J8").Value = 43122511.22151443
$ws.Range("K8").Value = 15096506.55049641
$ws.Range("M8").Value = 15096506.55049641
$ws.Range("N8").Value = 58219017.77201083
$ws.Range("O8").Value = 837860675.346328
$ws.Range("P8").Value = 820380499.322258
$ws.Range("Q8").Value = 0.01801791991760002
$ws.Range("R8").Value = 0.01840183495703287

# Row 9
$ws.Range("D9").Value = 9824
$ws.Range("E9").Value = 8726
$ws.Range("F9").Value = 0.8882328990228013
$ws.Range("G9").Value = 0.8851694055589369
$ws.Range("H9").Value = 0.09438221778644625
$ws.Range("I9").Value = 0.08354425161336274
$ws.Range("J9").Value = 45115835.6214844
$ws.Range("K9").Value = 15773246.61981758
$ws.Range("M9").Value = 15773246.61981758
$ws.Range("N9").Value = 60889082.24130198
$ws.Range("O9").Value = 874134762.184269
$ws.Range("P9").Value = 856685814.1783152
$ws.Range("Q9").Value = 0.01804441065860799
$ws.Range("R9").Value = 0.01841193861129403

# Row 10
$ws.Range("E10").Value = 8915
$ws.Range("F10").Value = 0.8889221258350782
$ws.Range("G10").Value = 0.8867117565148199
$ws.Range("H10").Value = 0.09325617978731912
$ws.Range("I10").Value = 0.08269135098507559
$ws.Range("J10").Value = 47276888.58374348
$ws.Range("K10").Value = 16520699.39275815
$ws.Range("M10").Value = 16520699.39275815
$ws.Range("N10").Value = 63797587.97650164
$ws.Range("O10").Value = 914424195.1217525
$ws.Range("P10").Value = 896938089.6587793
$ws.Range("Q10").Value = 0.01806677850486937
$ws.Range("R10").Value = 0.01841899634237085

# Row 11
$ws.Range("D11").Value = 10229
$ws.Range("E11").Value = 9106
$ws.Range("F11").Value = 0.8902140971746993
$ws.Range("G11").Value = 0.8880436902672127
$ws.Range("H11").Value = 0.09212643515120535
$ws.Range("I11").Value = 0.08181229944283946
$ws.Range("J11").Value = 49566607.15969561
$ws.Range("K11").Value = 17294514.02090722
$ws.Range("M11").Value = 17294514.02090722
$ws.Range("N11").Value = 66861121.18060283
$ws.Range("O11").Value = 955116216.3787864
$ws.Range("P11").Value = 937524390.505605
$ws.Range("Q11").Value = 0.01810723524983942
$ws.Range("R11").Value = 0.01844700169515624

# Row 12
$ws.Range("D12").Value = 9465
$ws.Range("E12").Value = 8384
$ws.Range("F12").Value = 0.8857897517168516
$ws.Range("G12").Value = 0.8845748048111416
$ws.Range("H12").Value = 0.09679007923257567
$ws.Range("I12").Value = 0.08561806544481057
$ws.Range("J12").Value = 41250664.96178105
$ws.Range("K12").Value = 14430500.42754652
$ws.Range("M12").Value = 14430500.42754652
$ws.Range("N12").Value = 55681165.38932757
$ws.Range("O12").Value = 801737041.7372
$ws.Range("P12").Value = 784037234.7332001
$ws.Range("Q12").Value = 0.0179990441707403
$ws.Range("R12").Value = 0.018405376413605

# Row 13
$ws.Range("D13").Value = 9648
$ws.Range("E13").Value = 8554
$ws.Range("F13").Value = 0.886608623548922
$ws.Range("G13").Value = 0.8849575832816057
$ws.Range("H13").Value = 0.101949910499804
$ws.Range("I13").Value = 0.09022134641168254
$ws.Range("J13").Value = 48011620.77939813
$ws.Range("K13").Value = 17541061.32943826
$ws.Range("M13").Value = 17541061.32943826
$ws.Range("N13").Value = 65552682.10883638
$ws.Range("O13").Value = 837547624.534428
$ws.Range("P13").Value = 820067448.510358
$ws.Range("Q13").Value = 0.02094335989453603
$ws.Range("R13").Value = 0.02138977880575723

# Row 14
$ws.Range("D14").Value = 9842
$ws.Range("E14").Value = 8741
$ws.Range("F14").Value = 0.8881324933956513
$ws.Range("G14").Value = 0.8866910123757354
$ws.Range("H14").Value = 0.1061872908635677
$ws.Range("I14").Value = 0.09415531643725347
$ws.Range("J14").Value = 54607314.52462393
$ws.Range("K14").Value = 20518986.07138735
$ws.Range("M14").Value = 20518986.07138735
$ws.Range("N14").Value = 75126300.5960113
$ws.Range("O14").Value = 875630888.769419
$ws.Range("P14").Value = 858181940.7634652
$ws.Range("Q14").Value = 0.02343337396448407
$ws.Range("R14").Value = 0.02390983204928902

# Row 15
$ws.Range("D15").Value = 10033
$ws.Range("E15").Value = 8884
$ws.Range("F15").Value = 0.8854779228545799
$ws.Range("G15").Value = 0.8836284066043366
$ws.Range("H15").Value = 0.1096290891937461
$ws.Range("I15").Value = 0.09687137740175458
$ws.Range("J15").Value = 60627454.10565276
$ws.Range("K15").Value = 23195982.15371279
$ws.Range("M15").Value = 23195982.15371279
$ws.Range("N15").Value = 83823436.25936554
$ws.Range("O15").Value = 914110715.330657
$ws.Range("P15").Value = 896624609.8676838
$ws.Range("Q15").Value = 0.02537546247373571
$ws.Range("R15").Value = 0.02587033848773776

# Row 16
$ws.Range("D16").Value = 10231
$ws.Range("E16").Value = 9096
$ws.Range("F16").Value = 0.889062652722119
$ws.Range("G16").Value = 0.8870684610883558
$ws.Range("H16").Value = 0.1086865141011976
$ws.Range("I16").Value = 0.09641237880480726
$ws.Range("J16").Value = 63928479.04605511
$ws.Range("K16").Value = 24475449.96408697
$ws.Range("M16").Value = 24475449.96408697
$ws.Range("N16").Value = 88403929.01014209
$ws.Range("O16").Value = 955558181.6979581
$ws.Range("P16").Value = 937966355.8247766
$ws.Range("Q16").Value = 0.0256137725916342
$ws.Range("R16").Value = 0.02609416618420722
